$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate for the data region so no stale styles / shared
# strings linger from the rows that are being removed (AdaBoost row, etc.)
$ws.Cells.Clear()

# ---- Header row ----
$ws.Range("A1").Value = "Model"
$ws.Range("B1").Value = "Parameter"
$ws.Range("C1").Value = "Accuracy"
$ws.Range("D1").Value = "F1"
$ws.Range("B1").HorizontalAlignment = -4108

# ---- Data rows ----
# Row 2: SVM
$ws.Range("A2").Value = "SVM"
$ws.Range("B2").Value = "C = 1, gamma = 10, kernel = rbf"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0

# Row 3: SVM
$ws.Range("A3").Value = "SVM"
$ws.Range("B3").Value = "C = 10, gamma = 10, kernel = rbf"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

# Row 4: SVM
$ws.Range("A4").Value = "SVM"
$ws.Range("B4").Value = "C = 1, gamma = 1, kernel = rbf"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0

# Row 5: Random Forest
$ws.Range("A5").Value = "Random Forest"
$ws.Range("B5").Value = "n_estimators=25, max_depth=10, min_samples_split=20, criterion='entropy'"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("B5").NumberFormat = "0.00E+00"

# Row 6: CNN
$ws.Range("A6").Value = "CNN"
$ws.Range("B6").Value = "solver = 'sgd', layers = 6, perceptrons = 300, epochs = 350, learning rate = 0.001"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0

# ---- Placeholder rows (reserved for pasted chart/diagram) ----
$ws.Range("B7").NumberFormat = "0.00E+00"
$ws.Range("B8").NumberFormat = "0.00E+00"
$ws.Range("B11").NumberFormat = "0.00E+00"
$ws.Range("B12").NumberFormat = "0.00E+00"
$ws.Range("B15").NumberFormat = "0.00E+00"
$ws.Range("B16").NumberFormat = "0.00E+00"
$ws.Range("B19").NumberFormat = "0.00E+00"
$ws.Range("B20").NumberFormat = "0.00E+00"

# ---- Column width for the Parameter column (narrower to fit the diagram) ----
$ws.Columns.Item(2).ColumnWidth = 67.36328125

# ---- Selection / view ----
$ws.Range("H20").Select()
